# Swap the order of slide 4 and slide 5 in the deck.
#
# Before: position 4 = sldId 327 ("2 SWOT 분석"), position 5 = sldId 325 ("구현 기능과 UI")
# After:  position 4 = sldId 325 ("구현 기능과 UI"), position 5 = sldId 327 ("2 SWOT 분석")
#
# i.e. the slide that used to sit at position 5 is moved to sit right before
# the slide that used to be at position 4 (equivalent to a drag-and-drop
# reorder of the two adjacent slides in the Slide pane).
$p = $ppt.ActivePresentation
$p.Slides.Item(5).MoveTo(4)
